$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while preserving it as literal text (avoids Excel
# auto-converting numeric-looking strings like "209.15" into numbers).
function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "28.291.60"
Set-TextValue "E2" "  -0.79%  "
# Row 3
Set-TextValue "D3" "1.552.32"
Set-TextValue "E3" "  -0.97%  "
# Row 4
Set-TextValue "E4" "  -0.10%  "
# Row 5
Set-TextValue "D5" "209.15"
Set-TextValue "E5" "  -1.49%  "
# Row 6
Set-TextValue "E6" "  -1.63%  "
# Row 7
Set-TextValue "E7" "  -0.06%  "
# Row 8
Set-TextValue "E8" "  -2.65%  "
# Row 9
Set-TextValue "E9" "  -2.23%  "
# Row 10
Set-TextValue "E10" "  -1.16%  "
# Row 11
Set-TextValue "E11" "  +0.04%  "
# Row 12
Set-TextValue "D12" "1.772.81"
Set-TextValue "E12" "  -1.09%  "
# Row 13
Set-TextValue "D13" "1.564.49"
Set-TextValue "E13" "  -0.26%  "
# Row 14
Set-TextValue "D14" "28.294.31"
Set-TextValue "E14" "  -0.74%  "
# Row 15
Set-TextValue "E15" "  -1.28%  "
# Row 16
Set-TextValue "E16" "  -2.30%  "
# Row 17
Set-TextValue "E17" "  -2.78%  "
# Row 18
Set-TextValue "D18" "226.56"
Set-TextValue "E18" "  -1.53%  "
# Row 19
Set-TextValue "E19" "  -0.55%  "
# Row 20
Set-TextValue "E20" "  -2.42%  "
# Row 21
Set-TextValue "E21" "  -0.07%  "
# Row 22
Set-TextValue "E22" "  +1.35%  "
# Row 23
Set-TextValue "E23" "  -3.10%  "
# Row 24
Set-TextValue "E24" "  -5.05%  "
# Row 25
Set-TextValue "D25" "147.81"
Set-TextValue "E25" "  -2.30%  "
# Row 26
Set-TextValue "D26" "14.77"
Set-TextValue "E26" "  -1.49%  "
# Row 27
Set-TextValue "E27" "  -0.13%  "
# Row 28
Set-TextValue "E28" "  -0.07%  "
# Row 29
Set-TextValue "E29" "  -3.06%  "
# Row 30
Set-TextValue "E30" "  -3.50%  "
# Row 31
Set-TextValue "E31" "  -4.48%  "
# Row 32
Set-TextValue "E32" "  -0.86%  "
# Row 33
Set-TextValue "E33" "  -0.97%  "
# Row 34
Set-TextValue "D34" "1.385.36"
Set-TextValue "E34" "  -0.68%  "
# Row 36
Set-TextValue "E36" "  -2.54%  "
# Row 37
Set-TextValue "E37" "  -1.37%  "
# Row 38
Set-TextValue "D38" "2.58"
Set-TextValue "E38" "  -1.36%  "
# Row 39
Set-TextValue "E39" "  -2.30%  "
# Row 40
Set-TextValue "D40" "1.93"
Set-TextValue "E40" "  +1.77%  "
# Row 41
Set-TextValue "D41" "0.512"
Set-TextValue "E41" "  -1.75%  "
# Row 42
Set-TextValue "E42" "  -0.04%  "
# Row 43
Set-TextValue "D43" "0.776"
Set-TextValue "E43" "  -1.37%  "
# Row 44
Set-TextValue "D44" "0.0466"
Set-TextValue "E44" "  +0.13%  "
# Row 45
Set-TextValue "E45" "  -0.83%  "
# Row 46
Set-TextValue "D46" "61.77"
# Row 47
Set-TextValue "B47" "RocketPoolETH"
Set-TextValue "C47" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D47" "1.687.00"
Set-TextValue "E47" "  -1.12%  "
# Row 48
Set-TextValue "B48" "WEMIXToken"
Set-TextValue "C48" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D48" "0.906"
Set-TextValue "E48" "  -6.53%  "
# Row 49
Set-TextValue "D49" "85.44"
Set-TextValue "E49" "  -0.98%  "
# Row 51
Set-TextValue "D51" "41.69"
Set-TextValue "E51" "  +5.17%  "
